$d = $word.ActiveDocument

# Step 1: "ตาราง " -> "ตารางที่ " (only the space run changes to "ที่ ")
$r1 = $d.Content
$r1.Find.Execute("ตาราง ", $true, $false, $false, $false, $false, $true, 1, $false, "ตารางที่ ", 2)

# Step 2: "… " -> "1 " (ellipsis + trailing space becomes "1" + space)
$r2 = $d.Content
$r2.Find.Execute("… ", $true, $false, $false, $false, $false, $true, 1, $false, "1 ", 2)
